$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 5 ("Austrittsgrund-Kategorie" / "Umsatzrückgang") is removed entirely;
# what used to be row 4 ("Austrittsgrund" / "betriebsbedingt") now carries the
# (renamed) reason value, so the sheet shrinks from 5 to 4 rows.
$ws.Rows.Item(5).Delete()

# Update the remaining values to the new dataset.
$ws.Range("B2").Value = "M100002"
$ws.Range("B3").Value = "31.12.2030"
$ws.Range("A4").Value = "Austrittsgrund"
$ws.Range("B4").Value = "Umsatzrueckgang"

# Drop the highlight fills that used to mark rows 2-4 - the refreshed sheet
# uses plain/default formatting for the data rows.
$ws.Range("A2:B4").Style = "Standard"

# The "value" column for Personalnummer/letzter Arbeitstag keeps a text number
# format so values like "M100002" / "31.12.2030" are not reinterpreted.
$ws.Range("B2:B3").NumberFormat = "@"

# Leave the sheet with B2:A2 (row with the employee number) selected, matching
# the saved view of the updated workbook.
$null = $ws.Range("A2:B2").Select()
